# Replace every "AxB=C" multiplication-answer cell in the quiz table with its
# updated value. wdReplaceAll (2) is used for values that occur exactly once
# in the document; the single value that appears twice ("45×98=4410") is
# replaced with wdReplaceOne (1) via two sequential calls so each occurrence
# (in document order) gets its own distinct replacement.
$d = $word.ActiveDocument

$d.Content.Find.Execute("26×85=2210", $true, $false, $false, $false, $false, $true, 1, $false, "75×56=4200", 2) | Out-Null
$d.Content.Find.Execute("14×26=364", $true, $false, $false, $false, $false, $true, 1, $false, "17×34=578", 2) | Out-Null
$d.Content.Find.Execute("16×70=1120", $true, $false, $false, $false, $false, $true, 1, $false, "56×62=3472", 2) | Out-Null
$d.Content.Find.Execute("22×66=1452", $true, $false, $false, $false, $false, $true, 1, $false, "66×85=5610", 2) | Out-Null
$d.Content.Find.Execute("18×22=396", $true, $false, $false, $false, $false, $true, 1, $false, "25×15=375", 2) | Out-Null
$d.Content.Find.Execute("89×25=2225", $true, $false, $false, $false, $false, $true, 1, $false, "55×79=4345", 2) | Out-Null
$d.Content.Find.Execute("99×63=6237", $true, $false, $false, $false, $false, $true, 1, $false, "50×36=1800", 2) | Out-Null
$d.Content.Find.Execute("18×23=414", $true, $false, $false, $false, $false, $true, 1, $false, "54×13=702", 2) | Out-Null
$d.Content.Find.Execute("96×37=3552", $true, $false, $false, $false, $false, $true, 1, $false, "53×48=2544", 2) | Out-Null
$d.Content.Find.Execute("86×21=1806", $true, $false, $false, $false, $false, $true, 1, $false, "37×19=703", 2) | Out-Null
$d.Content.Find.Execute("59×75=4425", $true, $false, $false, $false, $false, $true, 1, $false, "56×99=5544", 2) | Out-Null
$d.Content.Find.Execute("75×83=6225", $true, $false, $false, $false, $false, $true, 1, $false, "75×86=6450", 2) | Out-Null
$d.Content.Find.Execute("45×98=4410", $true, $false, $false, $false, $false, $true, 1, $false, "97×82=7954", 1) | Out-Null
$d.Content.Find.Execute("41×31=1271", $true, $false, $false, $false, $false, $true, 1, $false, "93×20=1860", 2) | Out-Null
$d.Content.Find.Execute("26×71=1846", $true, $false, $false, $false, $false, $true, 1, $false, "42×17=714", 2) | Out-Null
$d.Content.Find.Execute("64×19=1216", $true, $false, $false, $false, $false, $true, 1, $false, "43×36=1548", 2) | Out-Null
$d.Content.Find.Execute("69×50=3450", $true, $false, $false, $false, $false, $true, 1, $false, "17×90=1530", 2) | Out-Null
$d.Content.Find.Execute("52×40=2080", $true, $false, $false, $false, $false, $true, 1, $false, "52×81=4212", 2) | Out-Null
$d.Content.Find.Execute("17×43=731", $true, $false, $false, $false, $false, $true, 1, $false, "70×73=5110", 2) | Out-Null
$d.Content.Find.Execute("71×28=1988", $true, $false, $false, $false, $false, $true, 1, $false, "72×30=2160", 2) | Out-Null
$d.Content.Find.Execute("82×43=3526", $true, $false, $false, $false, $false, $true, 1, $false, "47×17=799", 2) | Out-Null
$d.Content.Find.Execute("56×49=2744", $true, $false, $false, $false, $false, $true, 1, $false, "40×93=3720", 2) | Out-Null
$d.Content.Find.Execute("100×49=4900", $true, $false, $false, $false, $false, $true, 1, $false, "48×68=3264", 2) | Out-Null
$d.Content.Find.Execute("27×67=1809", $true, $false, $false, $false, $false, $true, 1, $false, "46×90=4140", 2) | Out-Null
$d.Content.Find.Execute("39×69=2691", $true, $false, $false, $false, $false, $true, 1, $false, "45×61=2745", 2) | Out-Null
$d.Content.Find.Execute("46×40=1840", $true, $false, $false, $false, $false, $true, 1, $false, "92×41=3772", 2) | Out-Null
$d.Content.Find.Execute("22×82=1804", $true, $false, $false, $false, $false, $true, 1, $false, "14×52=728", 2) | Out-Null
$d.Content.Find.Execute("45×51=2295", $true, $false, $false, $false, $false, $true, 1, $false, "32×19=608", 2) | Out-Null
$d.Content.Find.Execute("33×13=429", $true, $false, $false, $false, $false, $true, 1, $false, "98×78=7644", 2) | Out-Null
$d.Content.Find.Execute("38×55=2090", $true, $false, $false, $false, $false, $true, 1, $false, "72×74=5328", 2) | Out-Null
$d.Content.Find.Execute("46×86=3956", $true, $false, $false, $false, $false, $true, 1, $false, "13×52=676", 2) | Out-Null
$d.Content.Find.Execute("16×88=1408", $true, $false, $false, $false, $false, $true, 1, $false, "44×19=836", 2) | Out-Null
$d.Content.Find.Execute("32×16=512", $true, $false, $false, $false, $false, $true, 1, $false, "38×46=1748", 2) | Out-Null
$d.Content.Find.Execute("100×67=6700", $true, $false, $false, $false, $false, $true, 1, $false, "90×27=2430", 2) | Out-Null
$d.Content.Find.Execute("82×37=3034", $true, $false, $false, $false, $false, $true, 1, $false, "56×80=4480", 2) | Out-Null
$d.Content.Find.Execute("51×16=816", $true, $false, $false, $false, $false, $true, 1, $false, "60×41=2460", 2) | Out-Null
$d.Content.Find.Execute("53×63=3339", $true, $false, $false, $false, $false, $true, 1, $false, "90×70=6300", 2) | Out-Null
$d.Content.Find.Execute("68×73=4964", $true, $false, $false, $false, $false, $true, 1, $false, "19×92=1748", 2) | Out-Null
$d.Content.Find.Execute("95×79=7505", $true, $false, $false, $false, $false, $true, 1, $false, "60×68=4080", 2) | Out-Null
$d.Content.Find.Execute("15×91=1365", $true, $false, $false, $false, $false, $true, 1, $false, "54×74=3996", 2) | Out-Null
$d.Content.Find.Execute("44×88=3872", $true, $false, $false, $false, $false, $true, 1, $false, "50×39=1950", 2) | Out-Null
$d.Content.Find.Execute("59×35=2065", $true, $false, $false, $false, $false, $true, 1, $false, "70×63=4410", 2) | Out-Null
$d.Content.Find.Execute("11×28=308", $true, $false, $false, $false, $false, $true, 1, $false, "19×61=1159", 2) | Out-Null
$d.Content.Find.Execute("44×75=3300", $true, $false, $false, $false, $false, $true, 1, $false, "22×75=1650", 2) | Out-Null
$d.Content.Find.Execute("73×31=2263", $true, $false, $false, $false, $false, $true, 1, $false, "99×37=3663", 2) | Out-Null
$d.Content.Find.Execute("91×94=8554", $true, $false, $false, $false, $false, $true, 1, $false, "71×54=3834", 2) | Out-Null
$d.Content.Find.Execute("79×78=6162", $true, $false, $false, $false, $false, $true, 1, $false, "15×11=165", 2) | Out-Null
$d.Content.Find.Execute("77×89=6853", $true, $false, $false, $false, $false, $true, 1, $false, "25×50=1250", 2) | Out-Null
$d.Content.Find.Execute("68×48=3264", $true, $false, $false, $false, $false, $true, 1, $false, "84×75=6300", 2) | Out-Null
$d.Content.Find.Execute("67×45=3015", $true, $false, $false, $false, $false, $true, 1, $false, "47×95=4465", 2) | Out-Null
$d.Content.Find.Execute("41×42=1722", $true, $false, $false, $false, $false, $true, 1, $false, "57×84=4788", 2) | Out-Null
$d.Content.Find.Execute("16×44=704", $true, $false, $false, $false, $false, $true, 1, $false, "16×50=800", 2) | Out-Null
$d.Content.Find.Execute("90×63=5670", $true, $false, $false, $false, $false, $true, 1, $false, "98×82=8036", 2) | Out-Null
$d.Content.Find.Execute("14×67=938", $true, $false, $false, $false, $false, $true, 1, $false, "94×53=4982", 2) | Out-Null
$d.Content.Find.Execute("93×84=7812", $true, $false, $false, $false, $false, $true, 1, $false, "55×20=1100", 2) | Out-Null
$d.Content.Find.Execute("66×16=1056", $true, $false, $false, $false, $false, $true, 1, $false, "70×35=2450", 2) | Out-Null
$d.Content.Find.Execute("31×92=2852", $true, $false, $false, $false, $false, $true, 1, $false, "25×71=1775", 2) | Out-Null
$d.Content.Find.Execute("49×45=2205", $true, $false, $false, $false, $false, $true, 1, $false, "36×41=1476", 2) | Out-Null
$d.Content.Find.Execute("44×16=704", $true, $false, $false, $false, $false, $true, 1, $false, "88×27=2376", 2) | Out-Null
$d.Content.Find.Execute("40×64=2560", $true, $false, $false, $false, $false, $true, 1, $false, "76×30=2280", 2) | Out-Null
$d.Content.Find.Execute("90×82=7380", $true, $false, $false, $false, $false, $true, 1, $false, "67×53=3551", 2) | Out-Null
$d.Content.Find.Execute("83×49=4067", $true, $false, $false, $false, $false, $true, 1, $false, "50×43=2150", 2) | Out-Null
$d.Content.Find.Execute("29×24=696", $true, $false, $false, $false, $false, $true, 1, $false, "94×91=8554", 2) | Out-Null
$d.Content.Find.Execute("81×12=972", $true, $false, $false, $false, $false, $true, 1, $false, "62×21=1302", 2) | Out-Null
$d.Content.Find.Execute("65×27=1755", $true, $false, $false, $false, $false, $true, 1, $false, "83×66=5478", 2) | Out-Null
$d.Content.Find.Execute("96×94=9024", $true, $false, $false, $false, $false, $true, 1, $false, "70×50=3500", 2) | Out-Null
$d.Content.Find.Execute("45×98=4410", $true, $false, $false, $false, $false, $true, 1, $false, "79×15=1185", 1) | Out-Null
$d.Content.Find.Execute("43×32=1376", $true, $false, $false, $false, $false, $true, 1, $false, "58×16=928", 2) | Out-Null
$d.Content.Find.Execute("83×40=3320", $true, $false, $false, $false, $false, $true, 1, $false, "48×60=2880", 2) | Out-Null
$d.Content.Find.Execute("35×59=2065", $true, $false, $false, $false, $false, $true, 1, $false, "55×72=3960", 2) | Out-Null
$d.Content.Find.Execute("67×85=5695", $true, $false, $false, $false, $false, $true, 1, $false, "62×52=3224", 2) | Out-Null
$d.Content.Find.Execute("55×76=4180", $true, $false, $false, $false, $false, $true, 1, $false, "10×34=340", 2) | Out-Null
$d.Content.Find.Execute("83×12=996", $true, $false, $false, $false, $false, $true, 1, $false, "69×66=4554", 2) | Out-Null
$d.Content.Find.Execute("13×31=403", $true, $false, $false, $false, $false, $true, 1, $false, "53×79=4187", 2) | Out-Null
$d.Content.Find.Execute("76×49=3724", $true, $false, $false, $false, $false, $true, 1, $false, "16×93=1488", 2) | Out-Null
$d.Content.Find.Execute("71×99=7029", $true, $false, $false, $false, $false, $true, 1, $false, "90×39=3510", 2) | Out-Null
$d.Content.Find.Execute("89×15=1335", $true, $false, $false, $false, $false, $true, 1, $false, "78×39=3042", 2) | Out-Null
$d.Content.Find.Execute("53×74=3922", $true, $false, $false, $false, $false, $true, 1, $false, "24×82=1968", 2) | Out-Null
$d.Content.Find.Execute("61×58=3538", $true, $false, $false, $false, $false, $true, 1, $false, "10×64=640", 2) | Out-Null
$d.Content.Find.Execute("26×21=546", $true, $false, $false, $false, $false, $true, 1, $false, "14×38=532", 2) | Out-Null
$d.Content.Find.Execute("37×83=3071", $true, $false, $false, $false, $false, $true, 1, $false, "45×89=4005", 2) | Out-Null
$d.Content.Find.Execute("81×20=1620", $true, $false, $false, $false, $false, $true, 1, $false, "87×40=3480", 2) | Out-Null
$d.Content.Find.Execute("62×22=1364", $true, $false, $false, $false, $false, $true, 1, $false, "40×25=1000", 2) | Out-Null
$d.Content.Find.Execute("61×22=1342", $true, $false, $false, $false, $false, $true, 1, $false, "46×82=3772", 2) | Out-Null
$d.Content.Find.Execute("85×13=1105", $true, $false, $false, $false, $false, $true, 1, $false, "63×41=2583", 2) | Out-Null
$d.Content.Find.Execute("83×18=1494", $true, $false, $false, $false, $false, $true, 1, $false, "91×55=5005", 2) | Out-Null
$d.Content.Find.Execute("69×28=1932", $true, $false, $false, $false, $false, $true, 1, $false, "32×23=736", 2) | Out-Null
$d.Content.Find.Execute("71×70=4970", $true, $false, $false, $false, $false, $true, 1, $false, "40×17=680", 2) | Out-Null
$d.Content.Find.Execute("44×78=3432", $true, $false, $false, $false, $false, $true, 1, $false, "55×46=2530", 2) | Out-Null
$d.Content.Find.Execute("93×10=930", $true, $false, $false, $false, $false, $true, 1, $false, "51×89=4539", 2) | Out-Null
$d.Content.Find.Execute("58×28=1624", $true, $false, $false, $false, $false, $true, 1, $false, "49×30=1470", 2) | Out-Null
$d.Content.Find.Execute("65×10=650", $true, $false, $false, $false, $false, $true, 1, $false, "66×18=1188", 2) | Out-Null
$d.Content.Find.Execute("25×51=1275", $true, $false, $false, $false, $false, $true, 1, $false, "41×47=1927", 2) | Out-Null
$d.Content.Find.Execute("77×31=2387", $true, $false, $false, $false, $false, $true, 1, $false, "59×64=3776", 2) | Out-Null
$d.Content.Find.Execute("85×27=2295", $true, $false, $false, $false, $false, $true, 1, $false, "76×15=1140", 2) | Out-Null
$d.Content.Find.Execute("32×83=2656", $true, $false, $false, $false, $false, $true, 1, $false, "51×53=2703", 2) | Out-Null
$d.Content.Find.Execute("27×89=2403", $true, $false, $false, $false, $false, $true, 1, $false, "47×76=3572", 2) | Out-Null
$d.Content.Find.Execute("59×56=3304", $true, $false, $false, $false, $false, $true, 1, $false, "44×48=2112", 2) | Out-Null
$d.Content.Find.Execute("53×10=530", $true, $false, $false, $false, $false, $true, 1, $false, "74×93=6882", 2) | Out-Null
$d.Content.Find.Execute("45×79=3555", $true, $false, $false, $false, $false, $true, 1, $false, "34×86=2924", 2) | Out-Null
